# Update cryptocurrency price (D) and volume-change (E) columns
# in place, preserving the original text (string) cell type for
# both columns so values like "1.00" / "40.699.22" are not
# auto-coerced into numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.699.22"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.380.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.98%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.06%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.68"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -6.38%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.48%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.23%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0823"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.96%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.07"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.77%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.84%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.744.76"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -4.10%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.11"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -4.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.362.85"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -4.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.759"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.625.13"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0911"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.31%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.08"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -3.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.77"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.11"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.65"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.89%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.60"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.91%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.39"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.60%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.92"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.66%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "155.16"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.99%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.23"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0731"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -5.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.41"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -6.96%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.84%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.02"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -8.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0994"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.12%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -8.44%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.81"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.26%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.90%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.962.50"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0270"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.55%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.73"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.33%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.86%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.606.15"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "93.67"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.94"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.62"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.92%  "
